$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

$ws.Range("A7").Value = "refrigerator_text"
$ws.Range("B7").Value = "text"
$ws.Rows.Item(7).RowHeight = 12.75

$ws.Range("B7").Select()
